$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 438, shifting existing rows 438:526 down to 439:527
$ws.Rows.Item(438).Insert()

# Populate the newly inserted row 438 with its data
$ws.Range("A438").Value = 8
$ws.Range("B438").Value = "Terminal La Palmera de La Serena"
$ws.Range("C438").Value = "Coquimbo"
$ws.Range("D438").Value = 44694
$ws.Range("E438").Value = 4
$ws.Range("F438").Value = 100112043
$ws.Range("G438").Value = "Pepino ensalada"
$ws.Range("H438").Value = "Sin especificar"
$ws.Range("I438").Value = "Primera"
$ws.Range("J438").Value = 800
$ws.Range("K438").Value = 18000
$ws.Range("L438").Value = 19000
$ws.Range("M438").Value = 18500
$ws.Range("N438").Value = "`$/caja 60 unidades"
$ws.Range("O438").Value = "Región de Arica y Parinacota"
$ws.Range("P438").Value = 308
$ws.Range("Q438").Value = 60
$ws.Range("R438").Value = "Hortaliza"
